$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on Price column cells we touch so numeric-looking
# strings (e.g. "1.00", "0.570") are preserved exactly as text, matching
# the original inline-string cell content instead of being coerced to numbers.
$priceCells = @('D2','D3','D5','D6','D8','D10','D13','D14','D15','D17','D18','D19','D20','D21','D24','D25','D26','D27','D29','D33','D36','D38','D39','D40','D41','D42','D43','D44','D46','D47','D49','D50','D51')
foreach ($ref in $priceCells) {
    $ws.Range($ref).NumberFormat = "@"
}

$ws.Range('D2').Value = '57.880.83'
$ws.Range('E2').Value = '  +0.26%  '
$ws.Range('D3').Value = '2.346.97'
$ws.Range('E3').Value = '  -0.59%  '
$ws.Range('E4').Value = '  +0.28%  '
$ws.Range('D5').Value = '540.68'
$ws.Range('E5').Value = '  -0.27%  '
$ws.Range('D6').Value = '134.50'
$ws.Range('E6').Value = '  -0.16%  '
$ws.Range('E7').Value = '  +0.22%  '
$ws.Range('D8').Value = '0.570'
$ws.Range('E8').Value = '  +6.31%  '
$ws.Range('E9').Value = '  +0.89%  '
$ws.Range('D10').Value = '5.54'
$ws.Range('E10').Value = '  +2.24%  '
$ws.Range('E11').Value = '  -1.56%  '
$ws.Range('E12').Value = '  +0.71%  '
$ws.Range('B13').Value = 'Avalanche'
$ws.Range('C13').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range('D13').Value = '23.77'
$ws.Range('E13').Value = '  +1.21%  '
$ws.Range('B14').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C14').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D14').Value = '2.767.37'
$ws.Range('E14').Value = '  -0.56%  '
$ws.Range('D15').Value = '57.804.06'
$ws.Range('E15').Value = '  +0.02%  '
$ws.Range('E16').Value = '  +0.68%  '
$ws.Range('D17').Value = '2.343.67'
$ws.Range('E17').Value = '  -0.37%  '
$ws.Range('D18').Value = '10.69'
$ws.Range('E18').Value = '  +1.11%  '
$ws.Range('B19').Value = 'BitcoinCash'
$ws.Range('C19').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D19').Value = '330.03'
$ws.Range('E19').Value = '  -2.42%  '
$ws.Range('B20').Value = 'Polkadot'
$ws.Range('C20').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D20').Value = '4.29'
$ws.Range('E20').Value = '  +1.42%  '
$ws.Range('D21').Value = '6.73'
$ws.Range('E21').Value = '  -1.67%  '
$ws.Range('E22').Value = '  +0.26%  '
$ws.Range('E23').Value = '  +1.16%  '
$ws.Range('D24').Value = '0.165'
$ws.Range('E24').Value = '  -2.53%  '
$ws.Range('D25').Value = '1.00'
$ws.Range('E25').Value = '  +0.08%  '
$ws.Range('D26').Value = '8.35'
$ws.Range('E26').Value = '  -1.66%  '
$ws.Range('D27').Value = '1.34'
$ws.Range('E27').Value = '  -5.87%  '
$ws.Range('E28').Value = '  +0.01%  '
$ws.Range('D29').Value = '170.17'
$ws.Range('E29').Value = '  -0.70%  '
$ws.Range('E30').Value = '  -0.22%  '
$ws.Range('E31').Value = '  -1.06%  '
$ws.Range('E32').Value = '  +0.91%  '
$ws.Range('D33').Value = '18.35'
$ws.Range('E33').Value = '  -1.09%  '
$ws.Range('E34').Value = '  -0.01%  '
$ws.Range('E35').Value = '  +0.27%  '
$ws.Range('D36').Value = '4.20'
$ws.Range('E36').Value = '  +1.45%  '
$ws.Range('E37').Value = '  -1.85%  '
$ws.Range('D38').Value = '1.60'
$ws.Range('E38').Value = '  -0.16%  '
$ws.Range('D39').Value = '39.06'
$ws.Range('D40').Value = '142.55'
$ws.Range('E40').Value = '  -4.26%  '
$ws.Range('D41').Value = '0.378'
$ws.Range('E41').Value = '  +0.06%  '
$ws.Range('D42').Value = '3.64'
$ws.Range('E42').Value = '  +0.38%  '
$ws.Range('D43').Value = '288.40'
$ws.Range('E43').Value = '  +0.92%  '
$ws.Range('D44').Value = '0.0949'
$ws.Range('E44').Value = '  +1.70%  '
$ws.Range('E45').Value = '  +0.40%  '
$ws.Range('D46').Value = '19.13'
$ws.Range('E46').Value = '  -0.25%  '
$ws.Range('D47').Value = '0.565'
$ws.Range('E47').Value = '  +0.81%  '
$ws.Range('E48').Value = '  +1.57%  '
$ws.Range('B49').Value = 'Polygon'
$ws.Range('C49').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D49').Value = '0.381'
$ws.Range('E49').Value = '  +0.43%  '
$ws.Range('B50').Value = 'EnergySwap'
$ws.Range('C50').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D50').Value = '17.45'
$ws.Range('E50').Value = '  -0.58%  '
$ws.Range('D51').Value = '11.08'
$ws.Range('E51').Value = '  +0.62%  '
